$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Haba" at Macroferia
# Regional de Talca. It belongs at the top of the existing (chronologically
# sorted) data block, so insert a fresh row at row 17 — this pushes the
# former rows 17-45 down to 18-46 untouched, matching the diff exactly.
$ws.Rows.Item(17).Insert()

$ws.Range('A17').Value = 5
$ws.Range('B17').Value = 'Macroferia Regional de Talca'
$ws.Range('C17').Value = 'Maule'
$ws.Range('D17').Value = 44495
$ws.Range('E17').Value = 7
$ws.Range('F17').Value = 100112026
$ws.Range('G17').Value = 'Haba'
$ws.Range('H17').Value = 'Sin especificar'
$ws.Range('I17').Value = 'Primera'
$ws.Range('J17').Value = 400
$ws.Range('K17').Value = 7000
$ws.Range('L17').Value = 7000
$ws.Range('M17').Value = 7000
$ws.Range('N17').Value = '$/saco 25 kilos'
$ws.Range('O17').Value = "Región del Maule"
$ws.Range('P17').Value = 280
$ws.Range('Q17').Value = 25
$ws.Range('R17').Value = 'Hortaliza'
